# Apply the "update new orleans xlsx files" edit:
#  1. hotel_info sheet gains a new "State" column (value "Louisiana") inserted
#     right after "Hotel_Name" (i.e. before "City"), shifting the remaining
#     columns one place to the right.
#  2. The two worksheets are reordered so that "review_info" becomes the
#     first (left-most) sheet and "hotel_info" becomes the second sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "State" column into hotel_info -----------------
$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Reorder the sheets: review_info first, hotel_info second ------
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wb.Worksheets.Item(1))
